$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters in the order the data appears (matching diff): B,C,D,E,F,I,J,K,L,M
$colIndexes = @(2,3,4,5,6,9,10,11,12,13)

$values = @(
    @(1.02,1.061376949907918,1.059268964046522,1.066137741294543,1.075198332236495,1.046764804841948,1.066352946779297,1.061998929533475,1.06884909798615,1.07788553521846),
    @(1.02,1.063115506579166,1.060615663001648,1.06770469937369,1.076849878206592,1.047272688990758,1.067741759823266,1.063158477468313,1.070229715825731,1.079352306830528),
    @(1.02,1.064237477795955,1.06148434080369,1.068715636053692,1.077915825873262,1.047598461532076,1.06863710703458,1.063905521738345,1.071119587538364,1.080298215184433),
    @(1.02,1.064708453617011,1.061848891014126,1.069139931483658,1.078363314474426,1.047734736153023,1.06901273142618,1.064218808669961,1.071492868959049,1.08069512389279),
    @(1.02,1.064787491869682,1.061910063266479,1.069211131755968,1.078438412891996,1.047757577524756,1.069075755054715,1.064271366036322,1.071555496791676,1.080761722931378),
    @(1.02,1.064243773731071,1.061489214448989,1.068721308251307,1.077921807717316,1.047600285105713,1.06864212919243,1.063909710907161,1.071124578550353,1.080303521632086),
    @(1.02,1.061965131824152,1.059724658319984,1.066667929341762,1.075757049959424,1.046937042460001,1.066822994961916,1.06239148549586,1.069316412467508,1.078381905973985),
    @(1.02,1.057926278518516,1.056593916734175,1.063026107187636,1.071921074080034,1.04574616961598,1.063591560669157,1.059690746258059,1.066102971527099,1.074970759497736),
    @(1.02,1.055216881918178,1.054491703971967,1.06058156944536,1.069348503321493,1.044937032602858,1.061419086660911,1.057872509297539,1.063941580305754,1.072679035291619),
    @(1.02,1.054039476293975,1.053577698911988,1.059518916211596,1.068230739508573,1.044582983561593,1.06047389372397,1.057080839238977,1.063000966382794,1.071682331557865),
    @(1.02,1.053601482877239,1.053237622533504,1.059123559119493,1.067814960944809,1.044450913589351,1.060122116153482,1.056786109889669,1.062650856330956,1.071311438857024),
    @(1.02,1.053695463783194,1.05331059627111,1.059208393803261,1.067904173888958,1.04447926850947,1.060197605143115,1.05684936076037,1.062725989145168,1.071391027289081),
    @(1.02,1.054003285042078,1.053549599883853,1.059486249030654,1.068196383283605,1.044572078085722,1.060444829861253,1.057056490516095,1.062972041039321,1.071651687262974),
    @(1.02,1.054192856847981,1.053696781408039,1.059657359453904,1.0683763442734,1.044629186677986,1.060597061171641,1.057184021120403,1.063123545308721,1.071812198882517),
    @(1.02,1.05529493197204,1.054552283690237,1.060652005350228,1.069422603629956,1.044960451437828,1.061481719959844,1.057924956870776,1.064003905052445,1.072745089785632),
    @(1.02,1.05598509190728,1.055087908847166,1.061274797149991,1.070077858936847,1.045167253356162,1.062035428334256,1.058388549771934,1.064554857667894,1.073329086053814),
    @(1.02,1.056387244719836,1.055399970663222,1.061637661563656,1.070459690487783,1.045287522075185,1.062357964164014,1.05865853570823,1.064875764911991,1.073669300429439),
    @(1.02,1.056524300167526,1.055506315140133,1.061761321685706,1.070589823303554,1.045328470520918,1.062467867575767,1.058750523081291,1.064985109431448,1.073785233842665),
    @(1.02,1.055911086386574,1.055030478594899,1.061208018926082,1.070007594406863,1.045145102268322,1.061976065551843,1.058338854170219,1.064495792747074,1.07326647241641),
    @(1.02,1.053912657518112,1.053479235235767,1.059404445402539,1.068110351359476,1.044544763500219,1.060372047536508,1.056995514503066,1.062899605004453,1.071574948090332),
    @(1.02,1.052652381141078,1.05250057807291,1.058266753853072,1.066914049699697,1.044164061188241,1.059359533188599,1.056147033275771,1.06189182087251,1.070507520495797),
    @(1.02,1.053320842545804,1.053019702792322,1.058870223216248,1.067548562611163,1.044366188483059,1.05989667089245,1.056597200323915,1.062426469389805,1.071073758942465),
    @(1.02,1.055944527549459,1.055056429963859,1.061238194365254,1.070039345059676,1.045155112489346,1.062002890375389,1.058361310771013,1.06452248304527,1.073294766123912),
    @(1.02,1.058973313626454,1.057405890747682,1.063970478271022,1.0729153901222,1.046056697567915,1.064430111371473,1.060392031867536,1.066937030238202,1.075855669620245)
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $values[$i]
    for ($j = 0; $j -lt $colIndexes.Count; $j++) {
        $ws.Cells.Item($rowNum, $colIndexes[$j]).Value = $rowData[$j]
    }
}

$wb.Save()
